$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the "Conclusion" paragraph that starts with "Overall, I found this
# assignment..." and append the new GitHub-link sentence + hyperlink to it,
# exactly as described by the commit ("Added GitHub link to document").
# ---------------------------------------------------------------------------

$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs.Item($i).Range.Text
    if ($ptext -like "*Overall, I found this assignment*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the 'Overall, I found this assignment' paragraph"
}

$para = $d.Paragraphs.Item($targetIndex)
$pr = $para.Range

# Insert the new sentence right before the paragraph's trailing mark (i.e. at
# the very end of the existing paragraph text).
$insertPoint = $d.Range($pr.End - 1, $pr.End - 1)

$newSentence = "I like how we had to add the code to GitHub because it will allow use to get familiar with source control software.  Here is my GitHub link for this assignment. "
$linkPlaceholder = "@@GITHUB_LINK_PLACEHOLDER@@"

$insertPoint.InsertAfter("  " + $newSentence + $linkPlaceholder)

# Find the placeholder we just inserted and turn it into a real hyperlink.
$linkUrl = "https://github.com/Slugdrew/Assignment05_python"
$findRange = $d.Content
$found = $findRange.Find.Execute($linkPlaceholder)
if (-not $found) {
    throw "Could not find hyperlink placeholder to replace"
}

$hlink = $d.Hyperlinks.Add($findRange, $linkUrl, $null, $null, $linkUrl)

# Trailing run after the hyperlink, matching the small Helvetica-styled
# whitespace run that Word leaves behind when pasting a rich-text hyperlink.
$para2 = $d.Paragraphs.Item($targetIndex)
$pr2 = $para2.Range
$tailInsertPoint = $d.Range($pr2.End - 1, $pr2.End - 1)
$tailInsertPoint.InsertAfter(" ")

$para3 = $d.Paragraphs.Item($targetIndex)
$pr3 = $para3.Range
$tailRange = $d.Range($pr3.End - 2, $pr3.End - 1)
$tailRange.Font.Name = "Helvetica"
$tailRange.Font.NameBi = "Helvetica"
$tailRange.Font.Size = 8.5
$tailRange.Font.Color = 2500134

Write-Output ("Updated paragraph " + $targetIndex + ": " + $d.Paragraphs.Item($targetIndex).Range.Text)
